# repair SOC by violation
#
# Adds a new "Greedy (rndstop!)" violation-repair result/solve-time/feasibility
# triple (columns B:D) for every Evs scenario, tags a couple of previously
# blank/ambiguous cells, and introduces a new column K ("minhstop!") flag.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Helper: write a numeric literal into a cell whose column style is
# "Text" (numFmtId 49) while still storing a genuine number (not a
# text string) - matches how the source file stores these figures.
# ------------------------------------------------------------------
function Set-NumericInTextCell($cell, $value) {
    $cell.NumberFormat = "General"
    $cell.Value = $value
    $cell.NumberFormat = "@"
}

# --- Header row: "Greedy" -> "Greedy(rndstop!)" ---
$ws.Range("H1").Value = "Greedy(rndstop!)"

# --- New column K ("minhstop!") ---
$ws.Columns("K").ColumnWidth = 8.25

$ws.Range("B1").Copy() | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null
$ws.Range("K1").Value = "minhstop!"

# --- New Greedy retry result / solve time / feasibility values ---
Set-NumericInTextCell $ws.Range("B4") 17925.951999409699
Set-NumericInTextCell $ws.Range("C4") 0.70399999618530196
$ws.Range("D4").Value = "Y"

Set-NumericInTextCell $ws.Range("B5") 34961.306080587899
Set-NumericInTextCell $ws.Range("C5") 2.8499999046325599
$ws.Range("D5").Value = "Y"

Set-NumericInTextCell $ws.Range("C6") 1.0909998416900599
$ws.Range("D6").Value = "N(failure)"

Set-NumericInTextCell $ws.Range("C7") 0.93099999427795399
$ws.Range("D7").Value = "N(failure)"

Set-NumericInTextCell $ws.Range("B9") 14352.025028387099
Set-NumericInTextCell $ws.Range("C9") 221.169999837875
$ws.Range("D9").Value = "Y"

Set-NumericInTextCell $ws.Range("B10") 10658.8759252645
Set-NumericInTextCell $ws.Range("C10") 24.8560001850128
$ws.Range("D10").Value = "Y"

$ws.Range("D10").Copy() | Out-Null
$ws.Range("K10").PasteSpecial(-4122) | Out-Null
$ws.Range("K10").Value = "Y"

Set-NumericInTextCell $ws.Range("B11") 7100.5354779756899
Set-NumericInTextCell $ws.Range("C11") 0.23100018501281699
$ws.Range("D11").Value = "Y"

Set-NumericInTextCell $ws.Range("B12") 3607.85144045987
Set-NumericInTextCell $ws.Range("C12") 12.317999839782701
$ws.Range("D12").Value = "Y"

# E12 previously held the (SOC-violating) computed total; once the violation
# is repaired via the new Greedy(rndstop!) columns this cell is cleared to a
# single space rather than a number.
$ws.Range("E12").Value = " "

# Leave the selection where the author left it
$ws.Range("I12").Select() | Out-Null
